# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures for rows 2-51 (row 49 is unchanged this run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.014.04'
$ws.Range("E2").Value = '  +1.52%  '

$ws.Range("D3").Value = '3.229.80'
$ws.Range("E3").Value = '  -1.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '3.224.62'
$ws.Range("E8").Value = '  -1.44%  '

$ws.Range("E9").Value = '  +1.17%  '

$ws.Range("E10").Value = '  -1.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.34'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.88%  '

$ws.Range("E12").Value = '  -0.60%  '

$ws.Range("E13").Value = '  +0.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.91%  '

$ws.Range("D15").Value = '3.765.68'
$ws.Range("E15").Value = '  -1.14%  '

$ws.Range("E16").Value = '  -1.85%  '

$ws.Range("D17").Value = '3.231.62'
$ws.Range("E17").Value = '  -1.20%  '

$ws.Range("D18").Value = '64.076.85'
$ws.Range("E18").Value = '  +1.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.710'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.98%  '

$ws.Range("E23").Value = '  -1.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  -0.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("E35").Value = '  -3.00%  '

$ws.Range("E36").Value = '  +0.73%  '

$ws.Range("D37").Value = '0.0₃0743'
$ws.Range("E37").Value = '  +3.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '51.67'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.19%  '

$ws.Range("E39").Value = '  +1.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '410.88'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.74%  '

$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("E43").Value = '  -2.44%  '

$ws.Range("D44").Value = '2.865.88'
$ws.Range("E44").Value = '  -7.06%  '

$ws.Range("E45").Value = '  +0.38%  '

$ws.Range("E46").Value = '  +1.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.62%  '

$ws.Range("E51").Value = '  +0.22%  '
